# Add the four new sheets required by the commit:
#   df1 Problem 12, df2 Problem 12, Problem 13, Sheet4
# and populate df1/df2/Problem 13 with their data.

$wb = $excel.ActiveWorkbook

# --- Sheet: df1 Problem 12 -------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$df1 = $wb.Worksheets.Add($null, $lastSheet)
$df1.Name = "df1 Problem 12"

$df1.Range("A1").Value = "order_id"
$df1.Range("B1").Value = "product"
$df1.Range("C1").Value = "quantity"

$df1.Range("A2").Value = 101
$df1.Range("B2").Value = "Apple"
$df1.Range("C2").Value = 5

$df1.Range("A3").Value = 102
$df1.Range("B3").Value = "Banana"
$df1.Range("C3").Value = 3

$df1.Range("A4").Value = 103
$df1.Range("B4").Value = "Kiwi"
$df1.Range("C4").Value = 2

$df1.Range("A1:C1").Select()

# --- Sheet: df2 Problem 12 -------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$df2 = $wb.Worksheets.Add($null, $lastSheet)
$df2.Name = "df2 Problem 12"

$df2.Range("A1").Value = "order_id"
$df2.Range("B1").Value = "product"
$df2.Range("C1").Value = "quantity"

$df2.Range("A2").Value = 201
$df2.Range("B2").Value = "Orange"
$df2.Range("C2").Value = 10

$df2.Range("A3").Value = 202
$df2.Range("B3").Value = "Papaya"
$df2.Range("C3").Value = 7

$df2.Range("A4").Value = 203
$df2.Range("B4").Value = "Mango"
$df2.Range("C4").Value = 8

$df2.Range("E6").Select()

# --- Sheet: Problem 13 ------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$p13 = $wb.Worksheets.Add($null, $lastSheet)
$p13.Name = "Problem 13"

$p13.Range("A1").Value = "Product"
$p13.Range("B1").Value = "Quarter"
$p13.Range("C1").Value = "Revenue"

$p13.Range("A2").Value = "Orange"
$p13.Range("B2").Value = "Q1"
$p13.Range("C2").Value = 5000

$p13.Range("A3").Value = "Orange"
$p13.Range("B3").Value = "Q2"
$p13.Range("C3").Value = 7000

$p13.Range("A4").Value = "Orange"
$p13.Range("B4").Value = "Q3"
$p13.Range("C4").Value = 8000

$p13.Range("A5").Value = "Orange"
$p13.Range("B5").Value = "Q4"
$p13.Range("C5").Value = 6000

$p13.Range("A6").Value = "Apple"
$p13.Range("B6").Value = "Q1"
$p13.Range("C6").Value = 3000

$p13.Range("A7").Value = "Apple"
$p13.Range("B7").Value = "Q2"
$p13.Range("C7").Value = 5000

$p13.Range("A8").Value = "Apple"
$p13.Range("B8").Value = "Q3"
$p13.Range("C8").Value = 4000

$p13.Range("A9").Value = "Apple"
$p13.Range("B9").Value = "Q4"
$p13.Range("C9").Value = 7000

$p13.Activate()
$p13.Range("F9").Select()
$excel.ActiveWindow.Zoom = 198

# --- Sheet: Sheet4 (blank) ---------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$sheet4 = $wb.Worksheets.Add($null, $lastSheet)
$sheet4.Name = "Sheet4"

# Re-activate Problem 13 so it remains the active/selected tab, matching the
# saved workbook state (activeTab points at "Problem 13").
$p13.Activate()
